$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Version and Date values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.2.0-ballot"
$meta.Range("B8").Value = "2025-12-19T09:47:21+00:00"

# --- Sheet "Include ValueSet #0": append version suffix to ValueSet URLs ---
$inc = $wb.Worksheets.Item("Include ValueSet #0")
$inc.Range("A2").Value = "https://mos.esante.gouv.fr/NOS/JDV_J283-PrestationsIndirects_SERAFIN/FHIR/JDV-J283-PrestationsIndirects-SERAFIN|20241025120000"
$inc.Range("A3").Value = "https://mos.esante.gouv.fr/NOS/JDV_J284-PrestationsDirects_SERAFIN/FHIR/JDV-J284-PrestationsDirects-SERAFIN|20241025120000"
